# Add trade credit columns: payment_term (Q) and tc_rate (R)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("Q1").Value = "payment_term"
$ws.Range("R1").Value = "tc_rate"

# Data rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 17).Value = 10
    $ws.Cells.Item($r, 18).Value = 0.15
}

# Autofit the new payment_term column to mirror the other best-fit columns
$ws.Columns.Item(17).AutoFit()

# Update the view selection to mimic the authored change
$ws.Range("R11").Select()
